$d = $word.ActiveDocument

# Manual line breaks (<w:br/>) show up in Range.Text as Chr(11).
$lb = [char]11

$line2053 = "LOT2053 -  Microbiologia  (Requisito fraco)"
$line2007 = "LOT2007 -  Bioquímica I  (Requisito fraco)"
$line2040 = "LOT2040 -  Engenharia Genética  (Requisito fraco)"

# Find the "Requisitos" bullet-list paragraph that holds the three
# LOT-requirement lines (one run per line, each ending in a manual break).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*$line2007*" -and $t -like "*$line2040*" -and $t -like "*$line2053*") {
        $target = $p
        break
    }
}

$r = $target.Range

# Insert a brand-new "LOT2053 ..." run + line break at the very start of the
# paragraph (this becomes its own run, just like typing it in Word would).
$insertPoint = $d.Range($r.Start, $r.Start)
$insertPoint.InsertBefore($line2053 + $lb)

# Re-fetch the paragraph (character offsets shifted after the insert above),
# then remove the now-duplicated trailing "LOT2053 ..." run.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*$line2007*" -and $t -like "*$line2040*" -and $t -like "*$line2053*") {
        $target = $p
        break
    }
}
$r2 = $target.Range
$full = $r2.Text
$needle = $line2053 + $lb
$lastIdx = $full.LastIndexOf($needle)

$delStart = $r2.Start + $lastIdx
$delEnd = $delStart + $needle.Length
$dup = $d.Range($delStart, $delEnd)
$dup.Delete()
